# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets to reflect that
# a.md has been handed back and is in sync with en-US: the Status column
# flips to "Handed back: in sync with en-US" (shared by both data rows),
# and the Latest Target File / Latest Handback File / Latest Handback
# DateTime columns (E, F, G) get populated for the a.md row (row 2) and
# the b.md row (row 3, which historically mirrors a.md's handoff file).

$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ Name = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackTime = "2016-01-26 09:00:40"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/041edb829a7b6d7c5e5b17a5e1d54fa97f69b19e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Name = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackTime = "2016-01-26 09:00:58"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b44af1cb5171fd8b2cde915c25331db78e43d6c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" }
)

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2ee2ce83501cb07379356a2474923e2b5142fcff/e2e/a.md"

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Name)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US" for both data rows.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Row 2 (a.md): Latest Target File, Latest Handback File, Latest Handback DateTime.
    $ws.Range("E2").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", "", "a.md")

    $ws.Range("F2").Value = $s.Xlf
    $ws.Hyperlinks.Add($ws.Range("F2"), $s.XlfUrl, "", "", $s.Xlf)

    $ws.Range("G2").Value = $s.HandbackTime

    # Row 3 (b.md): mirrors the same a.md target/handback file + datetime.
    $ws.Range("E3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, "", "", "a.md")

    $ws.Range("F3").Value = $s.Xlf
    $ws.Hyperlinks.Add($ws.Range("F3"), $s.XlfUrl, "", "", $s.Xlf)

    $ws.Range("G3").Value = $s.HandbackTime
}
